$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row at position 493, pushing the existing rows
# 493-528 down to 494-529 (dimension grows from A1:R528 to A1:R529).
$ws.Rows.Item(493).Insert()

# Populate the newly inserted row 493 with the new record.
$ws.Range("A493").Value2 = 3
$ws.Range("B493").Value2 = "Femacal de La Calera"
$ws.Range("C493").Value2 = "Coquimbo"
$ws.Range("D493").Value2 = 45021
$ws.Range("E493").Value2 = 5
$ws.Range("F493").Value2 = 100114013
$ws.Range("G493").Value2 = "Zanahoria"
$ws.Range("H493").Value2 = "Sin especificar"
$ws.Range("I493").Value2 = "Primera"
$ws.Range("J493").Value2 = 250
$ws.Range("K493").Value2 = 7500
$ws.Range("L493").Value2 = 8000
$ws.Range("M493").Value2 = 7760
$ws.Range("N493").Value2 = "$/saco 20 kilos"
$ws.Range("O493").Value2 = "Provincia de Quillota"
$ws.Range("P493").Value2 = 388
$ws.Range("Q493").Value2 = 20
$ws.Range("R493").Value2 = "Hortaliza"
